$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 33
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"
$ws.Range("H33").Value = "now()"

# Switch calculation to manual
$excel.Calculation = -4135

# Update view: scroll position / selection
$ws.Range("E31").Select() | Out-Null
